# Refresh the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") with the latest scraped values.
#
# Note: some new Price values (column D) look like plain numbers
# (e.g. "1.00", "407.88"). Cells in this sheet use the default General
# number format, so assigning those bare strings via .Value would make
# Excel auto-convert them to numeric cells (losing the trailing zero /
# the original text formatting). A leading apostrophe forces Excel to
# store them as literal text, exactly like typing '1.00 into the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.869.82"
$ws.Cells.Item(2, 5).Value = "  -0.52%  "
$ws.Cells.Item(3, 4).Value = "3.405.00"
$ws.Cells.Item(3, 5).Value = "  -0.71%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.11%  "
$ws.Cells.Item(5, 4).Value = "'407.88"
$ws.Cells.Item(5, 5).Value = "  -0.05%  "
$ws.Cells.Item(6, 4).Value = "'128.44"
$ws.Cells.Item(6, 5).Value = "  -3.99%  "
$ws.Cells.Item(7, 4).Value = "'0.634"
$ws.Cells.Item(7, 5).Value = "  +6.82%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
$ws.Cells.Item(9, 4).Value = "'0.725"
$ws.Cells.Item(9, 5).Value = "  +6.68%  "
$ws.Cells.Item(10, 4).Value = "'0.142"
$ws.Cells.Item(10, 5).Value = "  +16.20%  "
$ws.Cells.Item(11, 4).Value = "'42.26"
$ws.Cells.Item(11, 5).Value = "  -0.95%  "
$ws.Cells.Item(12, 4).Value = "'0.0000217"
$ws.Cells.Item(12, 5).Value = "  +64.81%  "
$ws.Cells.Item(13, 5).Value = "  -0.61%  "
$ws.Cells.Item(14, 4).Value = "3.950.27"
$ws.Cells.Item(14, 5).Value = "  -0.36%  "
$ws.Cells.Item(15, 4).Value = "'8.88"
$ws.Cells.Item(15, 5).Value = "  +5.28%  "
$ws.Cells.Item(16, 4).Value = "'20.79"
$ws.Cells.Item(16, 5).Value = "  +4.32%  "
$ws.Cells.Item(17, 4).Value = "3.391.63"
$ws.Cells.Item(17, 5).Value = "  -1.07%  "
$ws.Cells.Item(18, 4).Value = "'12.03"
$ws.Cells.Item(18, 5).Value = "  +9.03%  "
$ws.Cells.Item(19, 4).Value = "'1.06"
$ws.Cells.Item(19, 5).Value = "  +4.46%  "
$ws.Cells.Item(20, 4).Value = "61.911.93"
$ws.Cells.Item(20, 5).Value = "  -0.45%  "
$ws.Cells.Item(21, 4).Value = "'405.56"
$ws.Cells.Item(21, 5).Value = "  +28.58%  "
$ws.Cells.Item(22, 4).Value = "'89.08"
$ws.Cells.Item(22, 5).Value = "  +5.02%  "
$ws.Cells.Item(23, 4).Value = "'3.17"
$ws.Cells.Item(23, 5).Value = "  -1.48%  "
$ws.Cells.Item(24, 4).Value = "'13.05"
$ws.Cells.Item(24, 5).Value = "  +1.63%  "
$ws.Cells.Item(25, 5).Value = "  +3.12%  "
$ws.Cells.Item(26, 4).Value = "'32.69"
$ws.Cells.Item(26, 5).Value = "  +9.76%  "
$ws.Cells.Item(27, 4).Value = "'8.59"
$ws.Cells.Item(28, 4).Value = "'4.80"
$ws.Cells.Item(28, 5).Value = "  +0.21%  "
$ws.Cells.Item(29, 4).Value = "'7.59"
$ws.Cells.Item(29, 5).Value = "  -1.36%  "
$ws.Cells.Item(30, 5).Value = "  -0.95%  "
$ws.Cells.Item(31, 5).Value = "  +0.54%  "
$ws.Cells.Item(32, 5).Value = "  -2.09%  "
$ws.Cells.Item(33, 4).Value = "'11.80"
$ws.Cells.Item(33, 5).Value = "  +3.59%  "
$ws.Cells.Item(34, 4).Value = "'42.90"
$ws.Cells.Item(34, 5).Value = "  -0.50%  "
$ws.Cells.Item(35, 5).Value = "  +0.70%  "
$ws.Cells.Item(36, 4).Value = "'0.0493"
$ws.Cells.Item(36, 5).Value = "  +1.64%  "
$ws.Cells.Item(37, 4).Value = "'54.16"
$ws.Cells.Item(37, 5).Value = "  +3.73%  "
$ws.Cells.Item(38, 4).Value = "'0.998"
$ws.Cells.Item(38, 5).Value = "  -0.01%  "
$ws.Cells.Item(39, 5).Value = "  -2.75%  "
$ws.Cells.Item(40, 5).Value = "  +6.12%  "
$ws.Cells.Item(41, 4).Value = "'2.91"
$ws.Cells.Item(41, 5).Value = "  -2.96%  "
$ws.Cells.Item(42, 4).Value = "'0.310"
$ws.Cells.Item(42, 5).Value = "  +5.26%  "
$ws.Cells.Item(43, 4).Value = "'140.28"
$ws.Cells.Item(43, 5).Value = "  +1.90%  "
$ws.Cells.Item(44, 4).Value = "'1.96"
$ws.Cells.Item(44, 5).Value = "  -2.08%  "
$ws.Cells.Item(45, 4).Value = "'4.05"
$ws.Cells.Item(45, 5).Value = "  +0.63%  "
$ws.Cells.Item(46, 5).Value = "  +8.38%  "
$ws.Cells.Item(47, 4).Value = "'16.64"
$ws.Cells.Item(47, 5).Value = "  -1.26%  "
$ws.Cells.Item(48, 4).Value = "'21.75"
$ws.Cells.Item(48, 5).Value = "  +1.42%  "
$ws.Cells.Item(49, 4).Value = "2.113.63"
$ws.Cells.Item(49, 5).Value = "  -0.85%  "
$ws.Cells.Item(50, 5).Value = "  +4.42%  "
$ws.Cells.Item(51, 4).Value = "'0.130"
$ws.Cells.Item(51, 5).Value = "  +14.86%  "
